$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.67350212159806
$ws.Range("C2").Value = 8.81126097626035
$ws.Range("D2").Value = 14.51525121808
$ws.Range("E2").Value = 15.72606084443388
$ws.Range("G2").Value = 29.52208177275459
$ws.Range("H2").Value = 14.50513903438526
$ws.Range("I2").Value = 20.15252941985201
$ws.Range("J2").Value = 9.242276053008275
$ws.Range("M2").Value = 17.50286746444741
$ws.Range("O2").Value = 22.17773825202757
$ws.Range("B3").Value = 13.04582792188591
$ws.Range("C3").Value = 8.289227521766202
$ws.Range("D3").Value = 14.50879948820142
$ws.Range("E3").Value = 15.75826732128459
$ws.Range("G3").Value = 29.62996516628703
$ws.Range("H3").Value = 14.5671618011122
$ws.Range("I3").Value = 20.29641694033278
$ws.Range("J3").Value = 9.268780386819014
$ws.Range("M3").Value = 17.29728884361743
$ws.Range("O3").Value = 22.2794089944006
$ws.Range("B4").Value = 12.6446244193319
$ws.Range("C4").Value = 7.950714316067343
$ws.Range("D4").Value = 14.50797936462856
$ws.Range("E4").Value = 15.78118091454572
$ws.Range("G4").Value = 29.70838335612286
$ws.Range("H4").Value = 14.60815162224607
$ws.Range("I4").Value = 20.38975538163017
$ws.Range("J4").Value = 9.286032078911875
$ws.Range("M4").Value = 17.17182305337649
$ws.Range("O4").Value = 22.34786515515399
$ws.Range("B5").Value = 12.4773467639281
$ws.Range("C5").Value = 7.808308640843811
$ws.Range("D5").Value = 14.50843644136875
$ws.Range("E5").Value = 15.79130656441429
$ws.Range("G5").Value = 29.74338098205402
$ws.Range("H5").Value = 14.62558574328209
$ws.Range("I5").Value = 20.42904731462645
$ws.Range("J5").Value = 9.293308665519858
$ws.Range("M5").Value = 17.1209353591098
$ws.Range("O5").Value = 22.37727267096243
$ws.Range("B6").Value = 12.4493480616673
$ws.Range("C6").Value = 7.784395144455853
$ws.Range("D6").Value = 14.50856016139156
$ws.Range("E6").Value = 15.79303549052286
$ws.Range("G6").Value = 29.74937539028754
$ws.Range("H6").Value = 14.62852476210308
$ws.Range("I6").Value = 20.43564757862314
$ws.Range("J6").Value = 9.294531832089424
$ws.Range("M6").Value = 17.11250146834038
$ws.Range("O6").Value = 22.38224687856298
$ws.Range("B7").Value = 12.64238350129788
$ws.Range("C7").Value = 7.94881174194957
$ws.Range("D7").Value = 14.50798232356752
$ws.Range("E7").Value = 15.78131428283712
$ws.Range("G7").Value = 29.70884305921656
$ws.Range("H7").Value = 14.60838378832172
$ws.Range("I7").Value = 20.3902802008187
$ws.Range("J7").Value = 9.286129215281296
$ws.Range("M7").Value = 17.17113572474809
$ws.Range("O7").Value = 22.34825564385347
$ws.Range("B8").Value = 13.46046722865782
$ws.Range("C8").Value = 8.635029294615665
$ws.Range("D8").Value = 14.51237566284191
$ws.Range("E8").Value = 15.73651359796941
$ws.Range("G8").Value = 29.5567406145908
$ws.Range("H8").Value = 14.5259205913094
$ws.Range("I8").Value = 20.20110641286695
$ws.Range("J8").Value = 9.251212065103237
$ws.Range("M8").Value = 17.43185589356025
$ws.Range("O8").Value = 22.21153971315409
$ws.Range("B9").Value = 14.93226522774042
$ws.Range("C9").Value = 9.835958052161626
$ws.Range("D9").Value = 14.54583538462547
$ws.Range("E9").Value = 15.67360690164147
$ws.Range("G9").Value = 29.35594345712436
$ws.Range("H9").Value = 14.38731233867551
$ws.Range("I9").Value = 19.86970375585333
$ws.Range("J9").Value = 9.190477437762357
$ws.Range("M9").Value = 17.94698244356497
$ws.Range("O9").Value = 21.9915141638164
$ws.Range("B10").Value = 15.92511900178011
$ws.Range("C10").Value = 10.62807503636478
$ws.Range("D10").Value = 14.5854283956654
$ws.Range("E10").Value = 15.64264982952364
$ws.Range("G10").Value = 29.26890864732411
$ws.Range("H10").Value = 14.29960339142938
$ws.Range("I10").Value = 19.65029593426146
$ws.Range("J10").Value = 9.150542792987622
$ws.Range("M10").Value = 18.32489185979778
$ws.Range("O10").Value = 21.85947755525924
$ws.Range("B11").Value = 16.35636048701773
$ws.Range("C11").Value = 10.96860042113303
$ws.Range("D11").Value = 14.60666047585705
$ws.Range("E11").Value = 15.63188824767797
$ws.Range("G11").Value = 29.24264367659019
$ws.Range("H11").Value = 14.26277934117833
$ws.Range("I11").Value = 19.55570087665809
$ws.Range("J11").Value = 9.133386932911421
$ws.Range("M11").Value = 18.49608577466268
$ws.Range("O11").Value = 21.80590605163645
$ws.Range("B12").Value = 16.51664529809312
$ws.Range("C12").Value = 11.09468697807969
$ws.Range("D12").Value = 14.61515952813098
$ws.Range("E12").Value = 15.62829100035025
$ws.Range("G12").Value = 29.23462747551653
$ws.Range("H12").Value = 14.24927811984216
$ws.Range("I12").Value = 19.52062981803808
$ws.Range("J12").Value = 9.127035317377246
$ws.Range("M12").Value = 18.56075942310973
$ws.Range("O12").Value = 21.78655874646347
$ws.Range("B13").Value = 16.48226042838242
$ws.Range("C13").Value = 11.06765948612721
$ws.Range("D13").Value = 14.61330876741383
$ws.Range("E13").Value = 15.62904447172816
$ws.Range("G13").Value = 29.23626788978135
$ws.Range("H13").Value = 14.25216611920466
$ws.Range("I13").Value = 19.52814963588273
$ws.Range("J13").Value = 9.128396811446203
$ws.Range("M13").Value = 18.54683842398975
$ws.Range("O13").Value = 21.79068368467491
$ws.Range("B14").Value = 16.3696080973552
$ws.Range("C14").Value = 10.97903107977854
$ws.Range("D14").Value = 14.607350526615
$ws.Range("E14").Value = 15.63158272000804
$ws.Range("G14").Value = 29.24194543417101
$ws.Range("H14").Value = 14.2616596982429
$ws.Range("I14").Value = 19.55280052628485
$ws.Range("J14").Value = 9.132861479135759
$ws.Range("M14").Value = 18.50140988071587
$ws.Range("O14").Value = 21.80429548635792
$ws.Range("B15").Value = 16.30021021008644
$ws.Range("C15").Value = 10.9243704421808
$ws.Range("D15").Value = 14.60376056280328
$ws.Range("E15").Value = 15.63319971919056
$ws.Range("G15").Value = 29.24567476872663
$ws.Range("H15").Value = 14.26753254445291
$ws.Range("I15").Value = 19.56799759255551
$ws.Range("J15").Value = 9.135615080606749
$ws.Range("M15").Value = 18.47356208596494
$ws.Range("O15").Value = 21.81275556039751
$ws.Range("B16").Value = 15.89651829141163
$ws.Range("C16").Value = 10.605420909544
$ws.Range("D16").Value = 14.58410529218442
$ws.Range("E16").Value = 15.64341998697994
$ws.Range("G16").Value = 29.27089448935221
$ws.Range("H16").Value = 14.30207189889819
$ws.Range("I16").Value = 19.65658288831689
$ws.Range("J16").Value = 9.151684271499368
$ws.Range("M16").Value = 18.31368513328296
$ws.Range("O16").Value = 21.86310969399478
$ws.Range("B17").Value = 15.64357544960425
$ws.Range("C17").Value = 10.40466961758881
$ws.Range("D17").Value = 14.57286953285628
$ws.Range("E17").Value = 15.65054069402967
$ws.Range("G17").Value = 29.28978954381593
$ws.Range("H17").Value = 14.32404903900719
$ws.Range("I17").Value = 19.71226278513386
$ws.Range("J17").Value = 9.161800762601402
$ws.Range("M17").Value = 18.21538543875294
$ws.Range("O17").Value = 21.89566714012171
$ws.Range("B18").Value = 15.49617196759289
$ws.Range("C18").Value = 10.28733840177937
$ws.Range("D18").Value = 14.56671048558664
$ws.Range("E18").Value = 15.65494890811794
$ws.Range("G18").Value = 29.30191126039123
$ws.Range("H18").Value = 14.33697909713682
$ws.Range("I18").Value = 19.74477919988858
$ws.Range("J18").Value = 9.167714646847587
$ws.Range("M18").Value = 18.15878132734986
$ws.Range("O18").Value = 21.91500411985241
$ws.Range("B19").Value = 15.44593708369034
$ws.Range("C19").Value = 10.24729236904198
$ws.Range("D19").Value = 14.56467737973853
$ws.Range("E19").Value = 15.6564951235228
$ws.Range("G19").Value = 29.30623036753837
$ws.Range("H19").Value = 14.3414066639702
$ws.Range("I19").Value = 19.75587299306377
$ws.Range("J19").Value = 9.169733339833023
$ws.Range("M19").Value = 18.13960658525501
$ws.Range("O19").Value = 21.92165602662184
$ws.Range("B20").Value = 15.67070073260817
$ws.Range("C20").Value = 10.42623300908545
$ws.Range("D20").Value = 14.57403422043556
$ws.Range("E20").Value = 15.64975032925454
$ws.Range("G20").Value = 29.28764826905316
$ws.Range("H20").Value = 14.3216795768109
$ws.Range("I20").Value = 19.70628476707108
$ws.Range("J20").Value = 9.160714000416972
$ws.Range("M20").Value = 18.22585667679267
$ws.Range("O20").Value = 21.89213808936498
$ws.Range("B21").Value = 16.40277927865087
$ws.Range("C21").Value = 11.00514117958277
$ws.Range("D21").Value = 14.60908818766946
$ws.Range("E21").Value = 15.63082420193958
$ws.Range("G21").Value = 29.2402253320035
$ws.Range("H21").Value = 14.25885916600513
$ws.Range("I21").Value = 19.54553960446622
$ws.Range("J21").Value = 9.131546167535786
$ws.Range("M21").Value = 18.51475791388605
$ws.Range("O21").Value = 21.80027184174772
$ws.Range("B22").Value = 16.86362131888654
$ws.Range("C22").Value = 11.36680008431993
$ws.Range("D22").Value = 14.63467062077671
$ws.Range("E22").Value = 15.62124067285214
$ws.Range("G22").Value = 29.22048467606409
$ws.Range("H22").Value = 14.22038647132577
$ws.Range("I22").Value = 19.44485530214316
$ws.Range("J22").Value = 9.113327927706212
$ws.Range("M22").Value = 18.70265280383883
$ws.Range("O22").Value = 21.74570793155543
$ws.Range("B23").Value = 16.61929632036546
$ws.Range("C23").Value = 11.17530654163157
$ws.Range("D23").Value = 14.62077378182413
$ws.Range("E23").Value = 15.6261006042742
$ws.Range("G23").Value = 29.22998710342849
$ws.Range("H23").Value = 14.24068328304974
$ws.Range("I23").Value = 19.49819228655154
$ws.Range("J23").Value = 9.122974188246314
$ws.Range("M23").Value = 18.60247015341934
$ws.Range("O23").Value = 21.77432690201008
$ws.Range("B24").Value = 15.65844356328721
$ws.Range("C24").Value = 10.41649016858207
$ws.Range("D24").Value = 14.57350672871712
$ws.Range("E24").Value = 15.65010667379443
$ws.Range("G24").Value = 29.28861241961505
$ws.Range("H24").Value = 14.32274989148446
$ws.Range("I24").Value = 19.70898585503038
$ws.Range("J24").Value = 9.161205020991193
$ws.Range("M24").Value = 18.22112290750701
$ws.Range("O24").Value = 21.89373164460492
$ws.Range("B25").Value = 14.54916722230475
$ws.Range("C25").Value = 9.526810021621786
$ws.Range("D25").Value = 14.53413589754604
$ws.Range("E25").Value = 15.68794804464946
$ws.Range("G25").Value = 29.39971258936187
$ws.Range("H25").Value = 14.42233194689059
$ws.Range("I25").Value = 19.95512634510855
$ws.Range("J25").Value = 9.206082515844187
$ws.Range("M25").Value = 17.80751573780344
$ws.Range("O25").Value = 22.04585801213883
